$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.190.85'
$ws.Range('E2').Value = '  -3.68%  '
$ws.Range('D3').Value = '3.137.67'
$ws.Range('E3').Value = '  -5.23%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Formula = "'524.42"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -6.15%  '
$ws.Range('D6').Formula = "'135.23"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -5.15%  '
$ws.Range('D8').Value = '3.135.50'
$ws.Range('E8').Value = '  -5.30%  '
$ws.Range('D9').Formula = "'0.443"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -5.42%  '
$ws.Range('E10').Value = '  -7.38%  '
$ws.Range('D11').Formula = "'0.109"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -8.60%  '
$ws.Range('E12').Value = '  -6.37%  '
$ws.Range('D13').Value = '3.674.55'
$ws.Range('E13').Value = '  -5.26%  '
$ws.Range('E14').Value = '  -1.37%  '
$ws.Range('D15').Formula = "'25.60"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -4.79%  '
$ws.Range('D16').Value = '3.138.46'
$ws.Range('E16').Value = '  -5.18%  '
$ws.Range('D17').Value = '58.091.72'
$ws.Range('E17').Value = '  -3.87%  '
$ws.Range('D19').Formula = "'5.82"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -5.07%  '
$ws.Range('D20').Formula = "'13.11"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -8.14%  '
$ws.Range('D21').Formula = "'7.96"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -8.16%  '
$ws.Range('D22').Formula = "'345.46"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -7.78%  '
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('E24').Value = '  -8.02%  '
$ws.Range('D25').Formula = "'0.508"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -5.63%  '
$ws.Range('D26').Value = '3.265.66'
$ws.Range('E26').Value = '  -5.28%  '
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('D28').Value = '0.0₃0961'
$ws.Range('E28').Value = '  -6.54%  '
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('D30').Formula = "'6.81"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -5.41%  '
$ws.Range('D31').Formula = "'0.999"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -8.85%  '
$ws.Range('D33').Formula = "'6.88"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -9.49%  '
$ws.Range('D34').Formula = "'21.54"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.92%  '
$ws.Range('D35').Formula = "'1.23"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.20%  '
$ws.Range('D36').Formula = "'4.81"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -6.39%  '
$ws.Range('D37').Formula = "'157.14"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -5.58%  '
$ws.Range('D38').Formula = "'6.24"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -7.00%  '
$ws.Range('D39').Formula = "'1.37"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -10.65%  '
$ws.Range('D40').Formula = "'0.0694"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -5.01%  '
$ws.Range('D41').Value = '3.167.62'
$ws.Range('E41').Value = '  -5.20%  '
$ws.Range('D42').Formula = "'24.40"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -8.84%  '
$ws.Range('D43').Formula = "'40.53"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.40%  '
$ws.Range('D44').Formula = "'0.694"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -7.78%  '
$ws.Range('E45').Value = '  -2.62%  '
$ws.Range('D46').Formula = "'3.92"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -6.04%  '
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('E48').Value = '  -8.80%  '
$ws.Range('D49').Value = '2.265.18'
$ws.Range('E49').Value = '  -4.24%  '
$ws.Range('D50').Formula = "'6.20"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.39%  '
$ws.Range('D51').Formula = "'20.54"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.53%  '
